$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Time Period" row (row 6): end date 2023-04 -> 2023-05 ---
$ws.Range("B6").Value = "2008-12:2023-05"
$ws.Range("C6").Value = "2008-12:2023-05"
$ws.Range("D6").Value = "2010-09:2023-05"
$ws.Range("E6").Value = "2008-12:2023-05"
$ws.Range("F6").Value = "2009-12:2023-05"
$ws.Range("G6").Value = "2017-12:2023-05"
$ws.Range("H6").Value = "2017-12:2023-05"
$ws.Range("I6").Value = "2012-12:2023-05"
$ws.Range("J6").Value = "2012-12:2023-05"
$ws.Range("K6").Value = "2012-12:2023-05"
$ws.Range("L6").Value = "2012-12:2023-05"
$ws.Range("M6").Value = "2012-12:2023-05"
$ws.Range("N6").Value = "2020-06:2023-05"
$ws.Range("O6").Value = "2017-12:2023-05"
$ws.Range("P6").Value = "2017-12:2023-05"
$ws.Range("Q6").Value = "2018-02:2023-05"

# --- Update the "Update" row (row 8): refresh dates (typed as text, like the original) ---
$ws.Range("B8").Value = "'2023-06-21"
$ws.Range("C8").Value = "'2023-06-21"
$ws.Range("D8").Value = "'2023-06-21"
$ws.Range("E8").Value = "'2023-06-21"
$ws.Range("F8").Value = "'2023-06-21"
$ws.Range("G8").Value = "'2023-06-30"
$ws.Range("H8").Value = "'2023-06-30"
$ws.Range("I8").Value = "'2023-06-21"
$ws.Range("J8").Value = "'2023-06-21"
$ws.Range("K8").Value = "'2023-06-21"
$ws.Range("L8").Value = "'2023-06-21"
$ws.Range("M8").Value = "'2023-06-21"
$ws.Range("N8").Value = "'2023-06-30"
$ws.Range("O8").Value = "'2023-06-30"
$ws.Range("P8").Value = "'2023-06-30"
$ws.Range("Q8").Value = "'2023-06-29"

# --- Append a new data row (165) with the latest month's figures ---
$ws.Range("A164:Q164").Copy()
$ws.Range("A165:Q165").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(165, 1).Value = 45077
$newRowValues = @(434, 2202, 119, 1636, 6121, 1522, 340, 267240, 41700, 135095, 5676, 38260, 36881, 113679, 11871, 45392)
for ($i = 0; $i -lt $newRowValues.Length; $i++) {
    $ws.Cells.Item(165, 2 + $i).Value = $newRowValues[$i]
}
